$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet view: zoom level and selection cell ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 88
$ws.Range("O23").Select()

# --- Update "Team Member" column (C) for two tasks: Tyler -> Ben ---
$ws.Range("C17").Value = "Ben"
$ws.Range("C19").Value = "Ben"

# --- Populate weekly remaining-estimate values (Week 1..Week 4) ---
$values = @{
  4  = @(10, 0, 0, 0)
  5  = @(1, 1, 1, 0)
  6  = @(5, 5, 5, 0)
  7  = @(5, 5, 0, 0)
  8  = @(5, 5, 0, 0)
  9  = @(3, 0, 0, 0)
  10 = @(8, 8, 0, 0)
  11 = @(1, 0, 0, 0)
  12 = @(1, 1, 1, 0)
  13 = @(1, 0, 0, 0)
  14 = @(3, 3, 3, 0)
  15 = @(1, 0, 0, 0)
  16 = @(1, 1, 1, 0)
  17 = @(1, 1, 1, 0)
  18 = @(2, 2, 2, 0)
  19 = @(2, 2, 2, 0)
  20 = @(3, 3, 0, 0)
}

foreach ($row in $values.Keys) {
  $rowVals = $values[$row]
  $ws.Cells.Item($row, 5).Value = $rowVals[0]
  $ws.Cells.Item($row, 6).Value = $rowVals[1]
  $ws.Cells.Item($row, 7).Value = $rowVals[2]
  $ws.Cells.Item($row, 8).Value = $rowVals[3]
}

$wb.Application.Calculate()
